$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (Q0-Q3) with new values ---
$ws.Range("B2").Value = 0.2725968450730477
$ws.Range("C2").Value = 1.369443197887399
$ws.Range("D2").Value = 3.703882783350689
$ws.Range("E2").Value = 1.92454742299344
$ws.Range("F2").Value = 1.926195975934131
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.5701128739199846
$ws.Range("C3").Value = 1.3313120252481
$ws.Range("D3").Value = 2.906150528145031
$ws.Range("E3").Value = 1.704743537352476
$ws.Range("F3").Value = 1.624741117567991
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.1129289630813743
$ws.Range("C4").Value = 1.386823689668745
$ws.Range("D4").Value = 3.558153505617485
$ws.Range("E4").Value = 1.886306842912225
$ws.Range("F4").Value = 1.904692030384583
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.5895283849882742
$ws.Range("C5").Value = 1.252578073705098
$ws.Range("D5").Value = 2.867526754757556
$ws.Range("E5").Value = 1.693377322027656
$ws.Range("F5").Value = 1.606232434672746
$ws.Range("G5").Value = 43

# --- Add new rows 6-11 (Q4-Q9), copying the label style from A5 ---
$ws.Range("A5").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)

$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = 0.3216123161464985
$ws.Range("C6").Value = 1.332046485036206
$ws.Range("D6").Value = 2.890078795422161
$ws.Range("E6").Value = 1.7000231749662
$ws.Range("F6").Value = 1.689559483416096
$ws.Range("G6").Value = 42

$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.6210192535602802
$ws.Range("C7").Value = 1.275460222691501
$ws.Range("D7").Value = 2.981647775462735
$ws.Range("E7").Value = 1.726744849554425
$ws.Range("F7").Value = 1.631221148012869
$ws.Range("G7").Value = 41

$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.283533284100205
$ws.Range("C8").Value = 1.331940253175582
$ws.Range("D8").Value = 2.923668707468624
$ws.Range("E8").Value = 1.709873886422219
$ws.Range("F8").Value = 1.707683266217407
$ws.Range("G8").Value = 40

$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.7277015715620371
$ws.Range("C9").Value = 1.303080842420364
$ws.Range("D9").Value = 3.076077403806762
$ws.Range("E9").Value = 1.753874968122518
$ws.Range("F9").Value = 1.61664520434304
$ws.Range("G9").Value = 39

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.2602223904082975
$ws.Range("C10").Value = 1.380918912029693
$ws.Range("D10").Value = 3.046783774178897
$ws.Range("E10").Value = 1.745503874008562
$ws.Range("F10").Value = 1.749166497297726
$ws.Range("G10").Value = 38

$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.8291202883662652
$ws.Range("C11").Value = 1.24814772483617
$ws.Range("D11").Value = 2.865889069582671
$ws.Range("E11").Value = 1.692893697070986
$ws.Range("F11").Value = 1.496315835171674
$ws.Range("G11").Value = 37
